$d = $word.ActiveDocument

# "Text Box 22" (docPr id=22) was nudged/resized: position + extent change.
# NOTE: this runtime resolves $d.Shapes.Item(n) for *property writes* by raw
# XML document order (not by the z-order / relativeHeight order used when
# reading .Name while enumerating) -- index 2 in document order is the
# anchor that currently holds posOffset 1688934/3608871 and extent
# 289874x154001, i.e. "Text Box 22".
$shape = $d.Shapes.Item(2)

$shape.Left   = 1689652 / 12700.0
$shape.Top    = 3593990 / 12700.0
$shape.Height = 169572 / 12700.0
